$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1: reshuffle sheets.
#   before: Sheet1 (empty), clients, vendors, ...
#   after : clients, clients_structured, vendors, ...
# "Sheet1" becomes the new "clients_structured" sheet (same physical slot
# that used to be blank); "clients" moves in front of it.
# ---------------------------------------------------------------------------
$sheet1 = $wb.Worksheets.Item("Sheet1")
$sheet1.Name = "clients_structured"

$clients = $wb.Worksheets.Item("clients")
$structured = $wb.Worksheets.Item("clients_structured")
$clients.Move($structured)

# ---------------------------------------------------------------------------
# Step 2: trim the "clients" sheet - drop the extra 4th/5th Order/Quantity
# pair that used to live in H4:K4.
# ---------------------------------------------------------------------------
$clients = $wb.Worksheets.Item("clients")
$clients.Range("H4:K4").Clear()
$clients.Range("B14").Select()

$structured = $wb.Worksheets.Item("clients_structured")

# ---------------------------------------------------------------------------
# Step 3: build the normalised ("structured") client/order table.
# ---------------------------------------------------------------------------
$structured.Range("A1").Value = "Order No."
$structured.Range("B1").Value = "Name"
$structured.Range("C1").Value = "Email"
$structured.Range("D1").Value = "Order"
$structured.Range("E1").Value = "Quantity"

$data = @(
    @(1, "Carson Goble`n", "`nkcbku98@gmail.com", "Tires", 2),
    @(1, "Carson Goble`n", "`nkcbku98@gmail.com", "Mouse", 2),
    @(1, "Carson Goble`n", "`nkcbku98@gmail.com", "Tea", 1),
    @(2, "Aiden Herrera`n`n", "kcbradburn98@gmail.com`n`n", "Coffee", 1),
    @(2, "Aiden Herrera`n`n", "kcbradburn98@gmail.com`n`n", "Mouse", 1),
    @(2, "Aiden Herrera`n`n", "kcbradburn98@gmail.com`n`n", "Desktop", 2),
    @(3, "Cayden Doyle`n`n", "kenny.bradburn@revature.net`n", "Laptop", 1),
    @(3, "Cayden Doyle`n`n", "kenny.bradburn@revature.net`n", "Coffee", 1),
    @(3, "Cayden Doyle`n`n", "kenny.bradburn@revature.net`n", "Tires", 3)
)

$r = 2
foreach ($row in $data) {
    $structured.Cells.Item($r, 1).Value = $row[0]
    $structured.Cells.Item($r, 2).Value = $row[1]
    $structured.Cells.Item($r, 3).Value = $row[2]
    $structured.Cells.Item($r, 4).Value = $row[3]
    $structured.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}

$structured.Activate()
$structured.Range("A7").Select()
$structured.Range("C18").Select()

$wb.Worksheets | ForEach-Object { $_.Name }
